$p = $ppt.ActivePresentation
Write-Output ($p.SlideMaster.Name)
